# Update runtime_metrics_partition_0.xlsx: add rows for batch_size 5 and 25,
# and refresh the batch_size=1 measurements, per commit "json files for 1-25 runtime".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final table layout (row -> batch_size):
#   row2 -> 1   (updated measurements)
#   row3 -> 5   (new)
#   row4 -> 10  (previously row3, unchanged values)
#   row5 -> 25  (new)
#   row6 -> 100 (previously row4, unchanged values)

$data = @{
    2 = @(1,   6.456136703491211, 154.8913918534658, 23.3154914329185,  31.087321910558,   0.08320808410644531, 0.06937980651855469, 0.45013427734375,    0.2720355987548828, 12211.54141426086, 0.0667572021484375)
    3 = @(5,   25.94566345214844, 192.7104315224583, 29.0083158362126,  38.67775444828347, 0.1020431518554688,  0.0820159912109375,  1.123428344726562,   1.504182815551758,  80520.37000656128, 0.06842613220214844)
    4 = @(10,  50.34089088439941, 198.6456700371782, 29.90173541935637, 39.86898055914182, 0.1070499420166016,  0.08106231689453125, 1.799106597900391,   2.535343170166016,  257295.1803207397, 0.07867813110351562)
    5 = @(25,  123.5213279724121, 202.3941971024116, 30.46599370143181, 40.62132493524241, 0.1120567321777344,  0.08082389831542969, 3.931760787963867,   5.55109977722168,   1318080.536603928, 0.07724761962890625)
    6 = @(100, 489.4757270812988, 204.3002225999873, 30.75290390753089, 41.00387187670785, 0.1173019409179688,  0.08249282836914062, 14.51325416564941,   16.8766975402832,   17264854.43377495, 0.06580352783203125)
}

# Make sure every column-A cell (batch_size) ends up with the same bordered /
# centered / bold style already used by A1 (header) and A4 (existing data row),
# by copying that cell's formatting instead of re-deriving it property by
# property (which, for rows outside the original used range, produces a
# near-duplicate style entry instead of reusing the existing one).
$ws.Cells.Item(1, 1).Copy() | Out-Null
foreach ($r in 2..6) {
    $ws.Cells.Item($r, 1).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
}
$excel.CutCopyMode = 0

foreach ($r in 2..6) {
    $values = $data[$r]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = $i + 1
        $cell = $ws.Cells.Item($r, $col)
        $cell.Value = $values[$i]
    }
}
